# Applies the cryptos.xlsx price/volume refresh + the Quant/NEARProtocol row swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (from the commit diff).
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the original inlineStr cells) instead of auto-converting
# number-looking strings like "1.005" or "21.80" into floating point numbers.
$cellUpdates = [ordered]@{
    "D2" = "'28.768.83"
    "E2" = "'  +2.46%  "
    "D3" = "'1.876.30"
    "E3" = "'  +2.34%  "
    "D4" = "'1.005"
    "E4" = "'  +0.37%  "
    "D5" = "'325.11"
    "E5" = "'  +0.16%  "
    "E6" = "'  +0.22%  "
    "D7" = "'0.4588"
    "E7" = "'  -0.60%  "
    "D8" = "'0.3867"
    "E8" = "'  +0.03%  "
    "D9" = "'0.07853"
    "E9" = "'  +0.04%  "
    "D10" = "'0.9951"
    "E10" = "'  +3.55%  "
    "D11" = "'21.80"
    "D12" = "'1.883.51"
    "E12" = "'  +1.39%  "
    "D13" = "'6.994"
    "E13" = "'  +1.63%  "
    "D14" = "'5.715"
    "E14" = "'  +0.74%  "
    "D15" = "'0.06944"
    "E15" = "'  +1.19%  "
    "D16" = "'88.47"
    "E16" = "'  +0.20%  "
    "D18" = "'0.00001005"
    "E18" = "'  +1.15%  "
    "D19" = "'16.86"
    "E19" = "'  +1.02%  "
    "E20" = "'  +0.18%  "
    "D21" = "'28.777.90"
    "E21" = "'  +2.43%  "
    "D22" = "'5.277"
    "E22" = "'  -0.33%  "
    "D23" = "'11.04"
    "E23" = "'  +0.41%  "
    "D24" = "'2.129"
    "E24" = "'  +2.24%  "
    "D25" = "'2.106.91"
    "E25" = "'  +0.78%  "
    "D26" = "'153.26"
    "E26" = "'  -0.96%  "
    "D27" = "'19.26"
    "E27" = "'  +0.52%  "
    "D28" = "'5.788"
    "E28" = "'  +1.05%  "
    "D29" = "'1.970"
    "E29" = "'  +0.15%  "
    "D30" = "'119.01"
    "D31" = "'0.09318"
    "E31" = "'  +0.82%  "
    "D32" = "'0.9193"
    "E32" = "'  -2.41%  "
    "D33" = "'5.302"
    "E33" = "'  +0.77%  "
    "D34" = "'1.341"
    "E34" = "'  +1.48%  "
    "D35" = "'3.323"
    "E35" = "'  -0.13%  "
    "D36" = "'0.05769"
    "E36" = "'  -1.34%  "
    "D37" = "'1.155"
    "E37" = "'  +1.43%  "
    "D38" = "'0.02074"
    "E38" = "'  -1.68%  "
    "D39" = "'7.704"
    "E39" = "'  -0.28%  "
    "D40" = "'0.5645"
    "E40" = "'  +0.95%  "
    "D41" = "'0.1788"
    "E41" = "'  +1.56%  "
    "D42" = "'9.923"
    "E42" = "'  +0.24%  "
    "D43" = "'0.07218"
    "E43" = "'  -1.36%  "
    "D44" = "'11.77"
    "E44" = "'  +0.65%  "
    "D45" = "'0.5299"
    "E45" = "'  +0.45%  "
    "D46" = "'2.148"
    "E46" = "'  +1.89%  "
    "D47" = "'1.121"
    "E47" = "'  -1.83%  "
    "B48" = "'NEARProtocol"
    "C48" = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D48" = "'1.828"
    "E48" = "'  -0.20%  "
    "B49" = "'Quant"
    "C49" = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "D49" = "'113.60"
    "E49" = "'  +0.47%  "
    "D50" = "'2.413"
    "E50" = "'  +3.97%  "
    "D51" = "'1.003"
    "E51" = "'  +0.24%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = $cellUpdates[$ref]
    # Reset to the default style so no stray number-format/quote-prefix style
    # lingers on the cell (keeps formatting identical to the original file).
    $cell.Style = "Normal"
}
